$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill G (col 7) and H (col 8) values for rows 11-90:
#   G[row] = row - 1
#   H[row] = 100 - row
for ($r = 11; $r -le 90; $r++) {
    $g = $r - 1
    $h = 100 - $r
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Fill column I (col 9) with "=G+H" formulas. Applying the formula to a
# multi-cell range in one shot reproduces Excel's shared-formula grouping,
# matching how this block was originally built up (several fill operations).
$ws.Range("I11").Formula = "=G11+H11"
$ws.Range("I12").Formula = "=G12+H12"
$ws.Range("I13:I24").Formula = "=G13+H13"
$ws.Range("I25:I36").Formula = "=G25+H25"
$ws.Range("I37").Formula = "=G37+H37"
$ws.Range("I38").Formula = "=G38+H38"
$ws.Range("I39:I46").Formula = "=G39+H39"
$ws.Range("I47").Formula = "=G47+H47"
$ws.Range("I48").Formula = "=G48+H48"
$ws.Range("I49:I50").Formula = "=G49+H49"
$ws.Range("I51:I90").Formula = "=G51+H51"

# Update the view: select G11:G50 (active cell G50), scrolled so row 34 is
# at the top of the window.
$excel.ActiveWindow.ScrollRow = 34
[void]$ws.Range("G11:G50").Select()
